# Locate the paragraph holding the "{{发函日期}}" placeholder (the last
# paragraph of the letter body, right before the closing section break).
# Using Find keeps this content-addressed instead of relying on a
# hard-coded paragraph index / character offset.
$d = $word.ActiveDocument

$target = $d.Content
$found = $target.Find.Execute("{{发函日期}}", $false, $false, $false, $false, `
                               $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the {{发函日期}} placeholder paragraph"
}

# Grow the found range so it spans the whole paragraph, including its
# trailing paragraph mark (pilcrow) - this is the unit we are going to
# replace in one shot.
$target.Expand(4)  # wdParagraph

# Rebuild that paragraph (now carrying w:hint="eastAsia" + w:lang on its
# paragraph-mark run properties), followed by two brand new paragraphs:
#   1. a standalone paragraph holding just a manual page break, and
#   2. an (otherwise empty) paragraph that now owns the _GoBack bookmark
#      that used to sit at the end of the {{发函日期}} paragraph.
$newBodyXml = '<w:p w14:paraId="22944B1A"><w:pPr><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:pageBreakBefore w:val="0"/><w:widowControl w:val="0"/><w:kinsoku/><w:wordWrap/><w:overflowPunct/><w:topLinePunct w:val="0"/><w:autoSpaceDE/><w:autoSpaceDN/><w:bidi w:val="0"/><w:adjustRightInd/><w:snapToGrid/><w:spacing w:line="480" w:lineRule="auto"/><w:ind w:firstLine="4560" w:firstLineChars="1900"/><w:textAlignment w:val="auto"/><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="仿宋_GB2312" w:eastAsia="仿宋_GB2312"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="none"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="仿宋_GB2312" w:eastAsia="仿宋_GB2312"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="none"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>{{发函日期}}</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="仿宋_GB2312" w:eastAsia="仿宋_GB2312"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="none"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="仿宋_GB2312" w:eastAsia="仿宋_GB2312"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="none"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:br w:type="page"/></w:r></w:p><w:p><w:pPr><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:pageBreakBefore w:val="0"/><w:widowControl w:val="0"/><w:kinsoku/><w:wordWrap/><w:overflowPunct/><w:topLinePunct w:val="0"/><w:autoSpaceDE/><w:autoSpaceDN/><w:bidi w:val="0"/><w:adjustRightInd/><w:snapToGrid/><w:spacing w:line="480" w:lineRule="auto"/><w:ind w:firstLine="4560" w:firstLineChars="1900"/><w:textAlignment w:val="auto"/><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="仿宋_GB2312" w:eastAsia="仿宋_GB2312"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="none"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/></w:p>'

$package = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + $newBodyXml + '</w:body>' + `
  '</w:document>' + `
  '</pkg:xmlData>' + `
  '</pkg:part>' + `
  '</pkg:package>'

$target.InsertXML($package)

Write-Output "Paragraph count is now:"
Write-Output $d.Paragraphs.Count
